# Updates the cryptos price (D) and 1h-volume (E) columns to the latest
# scraped values, matching the GitHub Actions data-refresh commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that must stay plain TEXT even when it looks like a
# number (e.g. Price column values such as '0.675'). Going through a quoted
# formula and then Paste-Special-Values collapses it back to a literal text
# cell without touching the cell's number format / style.
function Set-TextValue($addr, $val) {
    $cell = $ws.Range($addr)
    $escaped = $val.Replace('"', '""')
    $cell.Formula = '="' + $escaped + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)  # xlPasteValues
}

$excel.CutCopyMode = 0

$ws.Range('D2').Value = '34.713.17'
$ws.Range('E2').Value = '  -2.07%  '
$ws.Range('D3').Value = '1.863.95'
$ws.Range('E3').Value = '  -3.02%  '
$ws.Range('E4').Value = '  -1.00%  '
Set-TextValue 'D5' '244.14'
$ws.Range('E5').Value = '  -3.92%  '
Set-TextValue 'D6' '0.675'
$ws.Range('E6').Value = '  -7.27%  '
$ws.Range('E7').Value = '  -0.99%  '
Set-TextValue 'D8' '41.14'
$ws.Range('E8').Value = '  +0.73%  '
$ws.Range('E9').Value = '  -4.40%  '
Set-TextValue 'D10' '51.07'
$ws.Range('E10').Value = '  -2.76%  '
Set-TextValue 'D11' '0.0725'
$ws.Range('E11').Value = '  -3.16%  '
$ws.Range('E12').Value = '  -3.09%  '
$ws.Range('D13').Value = '2.134.08'
$ws.Range('E13').Value = '  -3.16%  '
Set-TextValue 'D14' '12.69'
$ws.Range('E14').Value = '  -0.52%  '
$ws.Range('E15').Value = '  -2.36%  '
$ws.Range('D16').Value = '1.873.67'
$ws.Range('E16').Value = '  -2.67%  '
Set-TextValue 'D17' '4.80'
$ws.Range('E17').Value = '  -2.26%  '
$ws.Range('D18').Value = '34.703.97'
$ws.Range('E18').Value = '  -2.15%  '
Set-TextValue 'D19' '71.85'
$ws.Range('E19').Value = '  -3.52%  '
$ws.Range('D20').Value = '0.0₃0806'
$ws.Range('E20').Value = '  -3.83%  '
Set-TextValue 'D21' '240.95'
$ws.Range('E21').Value = '  -1.00%  '
Set-TextValue 'D22' '12.47'
$ws.Range('E22').Value = '  -4.56%  '
Set-TextValue 'D23' '4.83'
$ws.Range('E23').Value = '  -5.64%  '
$ws.Range('E24').Value = '  -0.86%  '
Set-TextValue 'D25' '2.46'
$ws.Range('E25').Value = '  +4.90%  '
$ws.Range('E26').Value = '  -10.41%  '
Set-TextValue 'D27' '163.58'
$ws.Range('E27').Value = '  -2.33%  '
$ws.Range('E28').Value = '  -3.99%  '
Set-TextValue 'D29' '17.99'
$ws.Range('E29').Value = '  -3.93%  '
$ws.Range('E30').Value = '  -7.06%  '
$ws.Range('D31').Value = '4.128.53'
$ws.Range('E31').Value = '  -0.06%  '
$ws.Range('E32').Value = '  +2.28%  '
Set-TextValue 'D33' '4.16'
$ws.Range('E33').Value = '  -5.17%  '
$ws.Range('E34').Value = '  -2.14%  '
$ws.Range('E35').Value = '  -0.91%  '
$ws.Range('E36').Value = '  -4.05%  '
Set-TextValue 'D37' '0.815'
$ws.Range('E37').Value = '  -11.30%  '
Set-TextValue 'D38' '1.58'
$ws.Range('E38').Value = '  -20.51%  '
$ws.Range('E39').Value = '  -4.82%  '
Set-TextValue 'D40' '96.29'
$ws.Range('E41').Value = '  +1.45%  '
$ws.Range('E42').Value = '  -4.06%  '
Set-TextValue 'D43' '0.0209'
$ws.Range('E43').Value = '  -0.89%  '
$ws.Range('E44').Value = '  -5.73%  '
$ws.Range('D45').Value = '1.272.68'
$ws.Range('E45').Value = '  -5.51%  '
Set-TextValue 'D46' '0.0809'
$ws.Range('E46').Value = '  +9.78%  '
Set-TextValue 'D47' '2.29'
$ws.Range('E47').Value = '  -6.51%  '
$ws.Range('E48').Value = '  -1.19%  '
$ws.Range('E49').Value = '  -2.80%  '
Set-TextValue 'D50' '11.82'
$ws.Range('E50').Value = '  -0.83%  '
Set-TextValue 'D51' '6.24'
$ws.Range('E51').Value = '  -8.21%  '

$excel.CutCopyMode = 0

